# Scheduler update: abbreviate month names on the "Routes" and "Shifts"
# sheets (e.g. "April, May, June" -> "Apr, May, Jun"), and move the active
# selection on the "Routes" and "Shifts" sheets.

$wb = $excel.ActiveWorkbook

# ---- Routes sheet: column D holds the "month" values ----
$ws = $wb.Worksheets.Item("Routes")

$ws.Range("D2:D6").Value = "Apr, May, Jun"
$ws.Range("D7:D11").Value = "Jul, Aug, Sep"
$ws.Range("D12:D16").Value = "Mar, Apr, May, Jun, Jul, Aug, Sep"
$ws.Range("D17:D20").Value = "Mar, Apr, Sep"
$ws.Range("D21:D24").Value = "May, Jun, Jul, Aug"
$ws.Range("D25:D26").Value = "Mar, Apr, May, Jun, Jul, Aug, Sep"

$ws.Range("B30").Select() | Out-Null

# ---- Shifts sheet: column D holds the "month" values ----
$ws = $wb.Worksheets.Item("Shifts")

$ws.Range("D2:D3").Value = "Mar, Apr, Sep"
$ws.Range("D4:D5").Value = "May, Jun, Jul, Aug"
$ws.Range("D6:D7").Value = "Mar, Apr, Sep"
$ws.Range("D8:D9").Value = "May, Jun, Jul, Aug"
$ws.Range("D10:D11").Value = "Mar, Apr, Sep"
$ws.Range("D12:D13").Value = "May, Jun, Jul, Aug"
$ws.Range("D14:D15").Value = "Mar, Apr, Sep"
$ws.Range("D16:D17").Value = "May, Jun, Jul, Aug"
$ws.Range("D18:D19").Value = "Mar"
$ws.Range("D20:D21").Value = "Apr"
$ws.Range("D22:D23").Value = "May, Jun, Jul"
$ws.Range("D24:D25").Value = "Aug"
$ws.Range("D26:D27").Value = "Sep"
$ws.Range("D28:D30").Value = "Mar, Apr, Sep"
$ws.Range("D31").Value = "May"
$ws.Range("D32:D33").Value = "Jun, Jul, Aug"
$ws.Range("D34").Value = "Mar, Apr, May, Jun, Jul, Aug, Sep"
$ws.Range("D35").Value = "Mar, Apr"
$ws.Range("D36").Value = "May, Jun, Jul, Aug, Sep"
$ws.Range("D37:D39").Value = "Mar, Apr, May, Jun, Jul, Aug, Sep"
$ws.Range("D40").Value = "Mar, Apr"
$ws.Range("D41").Value = "May, Jun, Jul, Aug, Sep"

$ws.Range("L23").Select() | Out-Null
